# Update reports - 2026-01-30 00:03
# Adds 5 new report rows to the "Reports" table (sheet1) and 1 new scan
# row to the "Scans" table (sheet2) describing the 2026-01-30 00:03:50 scan.

$wb = $excel.ActiveWorkbook
$wsReports = $wb.Worksheets.Item("Reports")
$wsScans   = $wb.Worksheets.Item("Scans")

# Helper: assign a value to a cell while forcing it to be stored as literal
# text (so digit-only / date-looking strings like "51305" or "2026-01-30"
# don't get silently reinterpreted by Excel as numbers or dates). The
# NumberFormat is reset back to the default ("Normal" style) immediately
# after the value is written so the cell's final formatting matches the
# rest of the sheet (no visible number-format override is left behind).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Reports table: append 5 new rows (28-32)
# ---------------------------------------------------------------------
$loReports = $wsReports.ListObjects.Item("Reports")

$newReports = @(
    @{ A="51305"; B="HC 571";   C="2024-26"; D="Women and Equalities Committee"; E="Commons"; F="Discrimination, harassment and abuse against Muslim women"; G="10th Report";        H="2026-01-30"; I="00:01:00"; K="0:02:50" },
    @{ A="51310"; B="HC 571";   C="2024-26"; D="Women and Equalities Committee"; E="Commons"; F="Large Print – 10th Report – Discrimination, harassment and abuse against Muslim women";               H="2026-01-30"; I="00:01:00"; K="0:02:50" },
    @{ A="51330"; B="HC 1639";  C="2024-26"; D="Justice Committee";              E="Commons"; F="Ending the cycle of reoffending – part one: rehabilitation in prisons: Government Response"; G="4th Special Report"; H="2026-01-30"; I="00:01:00"; K="0:02:50" },
    @{ A="51334"; B="HC 1234";  C="2024-26"; D="Public Accounts Committee";      E="Commons"; F="Costs of clinical negligence"; G="64th Report"; H="2026-01-30"; I="00:01:00"; K="0:02:50" },
    @{ A="51344"; B="HC 1651";  C="2024-26"; D="Scottish Affairs Committee";     E="Commons"; F="The work of the Committee in 2024-25, and Industrial transition in Scotland"; G="5th Report"; H="2026-01-30"; I="00:01:00"; K="0:02:50" }
)

foreach ($rowData in $newReports) {
    $newRow = $loReports.ListRows.Add()
    $r = $newRow.Range.Row

    Set-TextValue $wsReports.Cells.Item($r, 1) $rowData.A   # Publication ID
    Set-TextValue $wsReports.Cells.Item($r, 2) $rowData.B   # HC Number
    Set-TextValue $wsReports.Cells.Item($r, 3) $rowData.C   # Session
    Set-TextValue $wsReports.Cells.Item($r, 4) $rowData.D   # Committee Name
    Set-TextValue $wsReports.Cells.Item($r, 5) $rowData.E   # House
    Set-TextValue $wsReports.Cells.Item($r, 6) $rowData.F   # Report Title
    if ($rowData.ContainsKey("G")) {
        Set-TextValue $wsReports.Cells.Item($r, 7) $rowData.G   # Report Ordinal
    }
    Set-TextValue $wsReports.Cells.Item($r, 8) $rowData.H   # Publication Date
    Set-TextValue $wsReports.Cells.Item($r, 9) $rowData.I   # Publication Time
    Set-TextValue $wsReports.Cells.Item($r, 11) $rowData.K  # Late by max
}

# ---------------------------------------------------------------------
# Scans table: append 1 new row (15)
# ---------------------------------------------------------------------
$loScans = $wsScans.ListObjects.Item("Scans")
$newScanRow = $loScans.ListRows.Add()
$sr = $newScanRow.Range.Row

Set-TextValue $wsScans.Cells.Item($sr, 1) "2026-01-30"                               # Scan date
Set-TextValue $wsScans.Cells.Item($sr, 2) "00:03:50"                                 # Scan time
Set-TextValue $wsScans.Cells.Item($sr, 3) "51305, 51310, 51330, 51334, 51344"        # New Publication IDs
